$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "Sheet1" -> "OELogin"
$ws.Name = "OELogin"

# --- B2: "Pass@123" hyperlinked, with a thin border around it ---
$ws.Range("B2").Value = "Pass@123"
$ws.Hyperlinks.Add($ws.Range("B2"), "Pass@123") | Out-Null
$ws.Range("B2").Borders.LineStyle = 1

# --- A2: "110" stored as text (not a number), with a thin border ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("A2").Value = "110"

# --- Header row A1:B1 ("hrms_id" / "password"): bold, yellow fill, thin border ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.Color = 65535
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("A1").Value = "hrms_id"

$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Interior.Color = 65535
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("B1").Value = "password"

# Restore a sensible selection
$ws.Range("M19").Select() | Out-Null
